# Scheduled-runner refresh of cached market/profit figures across the
# per-job Leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Only the
# price/profit columns (H..N) on specific rows are refreshed; item/leve
# metadata columns (A..G) are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 6542.0835
$ws.Range("I80").Value = 509.77777
$ws.Range("J80").Value = 10161.467
$ws.Range("K80").Value = 1529.33331
$ws.Range("L80").Value = 30484.401
$ws.Range("M80").Value = -531.33331
$ws.Range("N80").Value = -32480.401

$ws.Range("H83").Value = 6542.0835
$ws.Range("I83").Value = 509.77777
$ws.Range("J83").Value = 10161.467
$ws.Range("K83").Value = 4587.99993
$ws.Range("L83").Value = 91453.20300000001
$ws.Range("M83").Value = 404.0000700000001
$ws.Range("N83").Value = -101437.203

$ws.Range("H86").Value = 6452.5454
$ws.Range("I86").Value = 2994.3333
$ws.Range("J86").Value = 10602.4
$ws.Range("K86").Value = 2994.3333
$ws.Range("L86").Value = 10602.4
$ws.Range("M86").Value = -1871.3333
$ws.Range("N86").Value = -12848.4

$ws.Range("H89").Value = 6452.5454
$ws.Range("I89").Value = 2994.3333
$ws.Range("J89").Value = 10602.4
$ws.Range("K89").Value = 14971.6665
$ws.Range("L89").Value = 53012
$ws.Range("M89").Value = -9355.666499999999
$ws.Range("N89").Value = -64244

$ws.Range("H98").Value = 1910.0769
$ws.Range("I98").Value = 1646.2122
$ws.Range("K98").Value = 1646.2122
$ws.Range("M98").Value = -148.2121999999999

$ws.Range("H106").Value = 3672.4375
$ws.Range("I106").Value = 3650.6
$ws.Range("K106").Value = 3650.6
$ws.Range("M106").Value = -3019.6

$ws.Range("H112").Value = 13984.667
$ws.Range("J112").Value = 20377.334
$ws.Range("L112").Value = 61132.00199999999
$ws.Range("N112").Value = -63348.00199999999

$ws.Range("H122").Value = 1910.0769
$ws.Range("I122").Value = 1646.2122
$ws.Range("K122").Value = 4938.6366
$ws.Range("M122").Value = -2488.6366

$ws.Range("H132").Value = 3503136
$ws.Range("I132").Value = 4487888
$ws.Range("J132").Value = 1794.3334
$ws.Range("K132").Value = 13463664
$ws.Range("L132").Value = 5383.0002
$ws.Range("M132").Value = -13461134
$ws.Range("N132").Value = -10443.0002

$ws.Range("H135").Value = 27669.37
$ws.Range("I135").Value = 911.4666999999999
$ws.Range("J135").Value = 128011.5
$ws.Range("K135").Value = 8203.2003
$ws.Range("L135").Value = 1152103.5
$ws.Range("M135").Value = -5668.2003
$ws.Range("N135").Value = -1157173.5

$ws.Range("H137").Value = 16428.791
$ws.Range("I137").Value = 30629.273
$ws.Range("J137").Value = 4413
$ws.Range("K137").Value = 91887.819
$ws.Range("L137").Value = 13239
$ws.Range("M137").Value = -89337.819
$ws.Range("N137").Value = -18339

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19550
$ws.Range("I32").Value = 21091.846
$ws.Range("K32").Value = 21091.846
$ws.Range("M32").Value = -20804.846

$ws.Range("H45").Value = 3108.5652
$ws.Range("I45").Value = 1893.9333
$ws.Range("K45").Value = 1893.9333
$ws.Range("M45").Value = -1516.9333

$ws.Range("H102").Value = 7542.857
$ws.Range("I102").Value = 8383.333000000001
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 8383.333000000001
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -6761.333000000001
$ws.Range("N102").Value = -5744

$ws.Range("H122").Value = 2111.889
$ws.Range("I122").Value = 2029.5454
$ws.Range("J122").Value = 2241.2856
$ws.Range("K122").Value = 6088.6362
$ws.Range("L122").Value = 6723.8568
$ws.Range("M122").Value = -3638.6362
$ws.Range("N122").Value = -11623.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 37999.8
$ws.Range("J93").Value = 37999.8
$ws.Range("L93").Value = 37999.8
$ws.Range("N93").Value = -41743.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2705282.5
$ws.Range("J31").Value = 4274.5
$ws.Range("L31").Value = 4274.5
$ws.Range("N31").Value = -4864.5

$ws.Range("H34").Value = 2705282.5
$ws.Range("J34").Value = 4274.5
$ws.Range("L34").Value = 4274.5
$ws.Range("N34").Value = -4678.5

$ws.Range("H62").Value = 6468.1113
$ws.Range("I62").Value = 5839.091
$ws.Range("J62").Value = 7456.5713
$ws.Range("K62").Value = 5839.091
$ws.Range("L62").Value = 7456.5713
$ws.Range("M62").Value = -5215.091
$ws.Range("N62").Value = -8704.5713

$ws.Range("H65").Value = 6468.1113
$ws.Range("I65").Value = 5839.091
$ws.Range("J65").Value = 7456.5713
$ws.Range("K65").Value = 29195.455
$ws.Range("L65").Value = 37282.85649999999
$ws.Range("M65").Value = -26075.455
$ws.Range("N65").Value = -43522.85649999999

$ws.Range("H69").Value = 72000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 72000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 72000
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -73498

$ws.Range("H72").Value = 72000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 72000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 216000
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -223488

$ws.Range("H94").Value = 2479.8572
$ws.Range("J94").Value = 2422
$ws.Range("L94").Value = 2422
$ws.Range("N94").Value = -3324

$ws.Range("H122").Value = 1929.7675
$ws.Range("I122").Value = 1962.2941
$ws.Range("K122").Value = 5886.8823
$ws.Range("M122").Value = -3436.8823

$ws.Range("H134").Value = 2447.8572
$ws.Range("I134").Value = 1827.1
$ws.Range("K134").Value = 5481.299999999999
$ws.Range("M134").Value = -2946.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40370396
$ws.Range("I4").Value = 45834656
$ws.Range("J4").Value = 23977610
$ws.Range("K4").Value = 137503968
$ws.Range("L4").Value = 71932830
$ws.Range("M4").Value = -137503856
$ws.Range("N4").Value = -71933054

$ws.Range("H25").Value = 2136.5557
$ws.Range("I25").Value = 2304.3333
$ws.Range("K25").Value = 6912.999899999999
$ws.Range("M25").Value = -6743.999899999999

$ws.Range("H30").Value = 2136.5557
$ws.Range("I30").Value = 2304.3333
$ws.Range("K30").Value = 6912.999899999999
$ws.Range("M30").Value = -6810.999899999999

$ws.Range("H37").Value = 43973.844
$ws.Range("J37").Value = 43973.844
$ws.Range("L37").Value = 131921.532
$ws.Range("N37").Value = -132145.532

$ws.Range("H98").Value = 1486.25
$ws.Range("I98").Value = 999
$ws.Range("J98").Value = 1648.6666
$ws.Range("K98").Value = 2997
$ws.Range("L98").Value = 4945.9998
$ws.Range("M98").Value = -1499
$ws.Range("N98").Value = -7941.9998

$ws.Range("H120").Value = 750
$ws.Range("I120").Value = 750
$ws.Range("K120").Value = 2250
$ws.Range("M120").Value = 2588

$ws.Range("H126").Value = 1934.3334
$ws.Range("I126").Value = 1934.3334
$ws.Range("K126").Value = 5803.0002
$ws.Range("M126").Value = -863.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 6598.65
$ws.Range("J24").Value = 8212.357
$ws.Range("L24").Value = 8212.357
$ws.Range("N24").Value = -8558.357

$ws.Range("H102").Value = 23038.621
$ws.Range("I102").Value = 32315.85
$ws.Range("K102").Value = 32315.85
$ws.Range("M102").Value = -30693.85

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H100").Value = 3154.3572
$ws.Range("I100").Value = 2417.1
$ws.Range("K100").Value = 2417.1
$ws.Range("M100").Value = -1876.1

$ws.Range("H106").Value = 14370
$ws.Range("J106").Value = 14370
$ws.Range("L106").Value = 14370
$ws.Range("N106").Value = -16894

$ws.Range("H132").Value = 3201.0952
$ws.Range("I132").Value = 3091.2632
$ws.Range("K132").Value = 9273.7896
$ws.Range("M132").Value = -6743.7896

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1121870
$ws.Range("J4").Value = 2515875
$ws.Range("L4").Value = 2515875
$ws.Range("N4").Value = -2516101

$ws.Range("H5").Value = 12873.875
$ws.Range("J5").Value = 12873.875
$ws.Range("L5").Value = 12873.875
$ws.Range("N5").Value = -13097.875

$ws.Range("H49").Value = 19850
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

$ws.Range("H62").Value = 4849.5
$ws.Range("J62").Value = 5499
$ws.Range("L62").Value = 5499
$ws.Range("N62").Value = -6747

$ws.Range("H65").Value = 4849.5
$ws.Range("J65").Value = 5499
$ws.Range("L65").Value = 27495
$ws.Range("N65").Value = -33735

$ws.Range("H81").Value = 16722.268
$ws.Range("I81").Value = 17488.143
$ws.Range("K81").Value = 34976.286
$ws.Range("M81").Value = -33915.286

$ws.Range("H84").Value = 16722.268
$ws.Range("I84").Value = 17488.143
$ws.Range("K84").Value = 174881.43
$ws.Range("M84").Value = -169577.43
